$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Extend header row formatting (style) to new columns F1:M1 to match existing header style
$ws.Range("B1").Copy()
$ws.Range("F1:M1").PasteSpecial(-4122)

$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"
$ws.Range("A2").Value = 43
$ws.Range("B2").Value = "台北富邦商業銀行金華分行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "蔣乃辛"
$ws.Range("F2").Value = 1433574
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-11-21"
$ws.Range("J2").Value = "蔣乃辛"
$ws.Range("K2").Value = 1722
$ws.Range("L2").Value = "tmp12421"
$ws.Range("M2").Value = 43
$ws.Range("A3").Value = 44
$ws.Range("B3").Value = "中華郵政股份有限公司青田支局"
$ws.Range("C3").Value = "活期儲蓄存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "蔣乃辛"
$ws.Range("F3").Value = 1419
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-11-21"
$ws.Range("J3").Value = "蔣乃辛"
$ws.Range("K3").Value = 1722
$ws.Range("L3").Value = "tmp12421"
$ws.Range("M3").Value = 44
$ws.Range("A4").Value = 45
$ws.Range("B4").Value = "f華郵政股份有限公司青ffl支局"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "楊際英"
$ws.Range("F4").Value = 132584
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2011-11-21"
$ws.Range("J4").Value = "蔣乃辛"
$ws.Range("K4").Value = 1722
$ws.Range("L4").Value = "tmp12421"
$ws.Range("M4").Value = 45
$ws.Range("A5").Value = 46
$ws.Range("B5").Value = "華南商業銀行信義分行"
$ws.Range("C5").Value = "活期儲蓄存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "楊際英"
$ws.Range("F5").Value = 1461520
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2011-11-21"
$ws.Range("J5").Value = "蔣乃辛"
$ws.Range("K5").Value = 1722
$ws.Range("L5").Value = "tmp12421"
$ws.Range("M5").Value = 46
$ws.Range("A6").Value = 47
$ws.Range("B6").Value = "臺灣中小企業銀行南京東"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "楊際英"
$ws.Range("F6").Value = 258909
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2011-11-21"
$ws.Range("J6").Value = "蔣乃辛"
$ws.Range("K6").Value = 1722
$ws.Range("L6").Value = "tmp12421"
$ws.Range("M6").Value = 47
$ws.Range("A7").Value = 48
$ws.Range("B7").Value = "臺灣中小企業銀行忠孝"
$ws.Range("C7").Value = "活期儲蓄存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "楊際英"
$ws.Range("F7").Value = 91020
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2011-11-21"
$ws.Range("J7").Value = "蔣乃辛"
$ws.Range("K7").Value = 1722
$ws.Range("L7").Value = "tmp12421"
$ws.Range("M7").Value = 48
$ws.Range("A8").Value = 49
$ws.Range("B8").Value = "永豐商業銀行三重分行"
$ws.Range("C8").Value = "活期儲蓄存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "楊際英"
$ws.Range("F8").Value = 68960
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2011-11-21"
$ws.Range("J8").Value = "蔣乃辛"
$ws.Range("K8").Value = 1722
$ws.Range("L8").Value = "tmp12421"
$ws.Range("M8").Value = 49
$ws.Range("A9").Value = 50
$ws.Range("B9").Value = "台北富邦商業銀行金華分行"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "美金"
$ws.Range("E9").Value = "楊際英"
$ws.Range("F9").Value = 24
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2011-11-21"
$ws.Range("J9").Value = "蔣乃辛"
$ws.Range("K9").Value = 1722
$ws.Range("L9").Value = "tmp12421"
$ws.Range("M9").Value = 50
$ws.Range("A10").Value = 51
$ws.Range("B10").Value = "台北富邦商業銀行金華分行"
$ws.Range("C10").Value = "活期存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "楊際英"
$ws.Range("F10").Value = 39854
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2011-11-21"
$ws.Range("J10").Value = "蔣乃辛"
$ws.Range("K10").Value = 1722
$ws.Range("L10").Value = "tmp12421"
$ws.Range("M10").Value = 51
